# Apply the commit "Update gh-pages to output generated at 456a3b4"
# This updates attendance / price figures scraped from bilibili show listings
# across the 4 worksheets: 展览 (sheet1), 演出 (sheet2), 本地生活 (sheet3),
# 全部类型 (sheet4). One row (演出 row 6 / original spreadsheet row 5) also
# had its event details replaced wholesale because the event was cancelled.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) - column F ("想去人数" / interest count) bumps
# ---------------------------------------------------------------------------
$ws1.Range("F5").Value  = 1951
$ws1.Range("F6").Value  = 1951
$ws1.Range("F7").Value  = 1315
$ws1.Range("F13").Value = 1753
$ws1.Range("F14").Value = 5
$ws1.Range("F15").Value = 1869
$ws1.Range("F17").Value = 1035
$ws1.Range("F19").Value = 523
$ws1.Range("F20").Value = 1606
$ws1.Range("F26").Value = 448
$ws1.Range("F28").Value = 1030
$ws1.Range("F29").Value = 4579
$ws1.Range("F31").Value = 36

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------------
# Row 5 (原"上海·2024·松田瑠华专场见面会") -> event changed to a cancelled show
$ws2.Range("C5").Value = "上海·【漫乐季】《冠位时之门：热血番同人Only》Fun肆二次元·动漫ACG超燃音乐演唱会（取消）"
$ws2.Range("D5").Value = "嘉定区城中路149号 嘉定影剧院"
$ws2.Range("E5").Value = "2024.10.27 15:00-10.27 17:30"
$ws2.Range("F5").Value = 48
$ws2.Range("G5").Value = "不可售"
$ws2.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92750"
$ws2.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202409/Jud6aOcf1727073479811.jpeg"

$ws2.Range("F6").Value  = 29
$ws2.Range("F8").Value  = 171
$ws2.Range("F16").Value = 30
$ws2.Range("F18").Value = 1
$ws2.Range("F33").Value = 476
$ws2.Range("F40").Value = 46
$ws2.Range("F43").Value = 100

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life)
# ---------------------------------------------------------------------------
$ws3.Range("F4").Value  = 9597
$ws3.Range("F5").Value  = 174
$ws3.Range("F9").Value  = 3104
$ws3.Range("F10").Value = 619
$ws3.Range("F11").Value = 893
$ws3.Range("F13").Value = 39
$ws3.Range("F14").Value = 58
$ws3.Range("F15").Value = 10
$ws3.Range("F16").Value = 313

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types) - aggregated view, mirrors the above updates
# ---------------------------------------------------------------------------
$ws4.Range("F3").Value  = 174
$ws4.Range("F9").Value  = 3104
$ws4.Range("F10").Value = 619
$ws4.Range("F11").Value = 893
$ws4.Range("F12").Value = 1951
$ws4.Range("F13").Value = 39
$ws4.Range("F14").Value = 39
$ws4.Range("F15").Value = 58
$ws4.Range("F16").Value = 1315
$ws4.Range("F17").Value = 10
$ws4.Range("F18").Value = 1753
$ws4.Range("F19").Value = 5
$ws4.Range("F22").Value = 30
$ws4.Range("F23").Value = 1869
$ws4.Range("F24").Value = 1035
$ws4.Range("F26").Value = 523
$ws4.Range("F27").Value = 1606
$ws4.Range("F28").Value = 1
$ws4.Range("F35").Value = 448
$ws4.Range("F37").Value = 1030
$ws4.Range("F39").Value = 313
$ws4.Range("F41").Value = 4579
$ws4.Range("F42").Value = 476
$ws4.Range("F43").Value = 36
$ws4.Range("F45").Value = 46
$ws4.Range("F47").Value = 100
